# Update the multiplication problems in the document to match the new
# randomly-generated worksheet values.
#
# Each "old" text is unique within the document, so a simple
# Find/Replace (wdReplaceAll) for each pair is sufficient and safe.

$d = $word.ActiveDocument

$replacements = @(
    @("40×91=", "13×18="),
    @("81×12=", "30×94="),
    @("69×42=", "52×34="),
    @("85×30=", "21×80="),
    @("42×65=", "79×53="),
    @("20×42=", "32×42="),
    @("65×59=", "58×50="),
    @("93×85=", "35×87="),
    @("86×47=", "86×41="),
    @("18×64=", "35×65="),
    @("63×95=", "73×29="),
    @("70×87=", "12×98="),
    @("14×51=", "53×27="),
    @("85×63=", "40×17="),
    @("94×84=", "50×96="),
    @("53×65=", "69×40="),
    @("88×80=", "44×85="),
    @("50×77=", "88×40="),
    @("14×84=", "47×71="),
    @("83×44=", "23×90="),
    @("36×74=", "97×85="),
    @("23×73=", "36×73="),
    @("42×53=", "93×53="),
    @("48×49=", "43×21="),
    @("69×20=", "64×91=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}
